$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3 and F27 currently hold a stray numeric 0; the rest of their rows use the
# text "/" (shared string) for "no data" cells, so bring these two in line
# with the surrounding cells (C3/E3 and C27/E27).
$ws.Range("F3").Value = "/"
$ws.Range("F27").Value = "/"

# Move the selection/active cell to M9 (and drop the old frozen scroll
# position at A16 that was previously saved in the view).
$ws.Range("M9").Select() | Out-Null
